$d = $word.ActiveDocument

# 1. The empty paragraph right after the second screenshot (before "In Chrome I get an error:")
#    gets new commentary text added to it.
$p = $d.Paragraphs(53)
$p.Range.InsertBefore("It would be nice if the tutorial said that Chrome/FF weren" + [char]0x2019 + "t expected to work yet.  Not telling me would cause me to do a bunch of needless debugging because there is no point in continuing a tutorial if something breaks in the middle.")

# 2. The empty paragraph right after the "Now the tutorial is discussing..." paragraph
#    (and right before the bookmark paragraph) is removed.
$p2 = $d.Paragraphs(55)
$p2.Range.Delete()

# 3. A new empty paragraph is inserted right after the bookmark paragraph.
$bm = $d.Paragraphs(55)
$insertionPoint = $d.Range($bm.Range.End, $bm.Range.End)
$null = $insertionPoint.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>')
